$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    'coChannelProfile_1.0.0-tsp.d.t+gendoc.${date}.${time}docx''',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'CoChannelProfile_1.0.0-tsi.d.t+gendoc.1.docx''',
    2)

Write-Output "Found1: $found"
